# Adds SAM0-UNICON (SNU) citation rows to the "Citations" sheet of the
# ES-DOC CMIP6 model-citations workbook and makes that sheet the active tab.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Citations")

# ---------------------------------------------------------------------
# BibTeX blocks (kept as variables to keep the per-cell code readable).
# NOTE: each here-string is prefixed with a single leading space because
# a line that starts with the two characters  "@  is (mis-)parsed as the
# terminator of a here-string by this runtime; TrimStart() removes it.
# ---------------------------------------------------------------------

$bibJCLI = @"
 "@article{doi:10.1175/JCLI-D-18-0796.1,
author = {Park, Sungsu and Shin, Jihoon and Kim, Siyun and Oh, Eunsil and Kim, Yoonjae},
title = {Global Climate Simulated by the Seoul National University Atmosphere Model Version 0 with a Unified Convection Scheme (SAM0-UNICON)},
journal = {Journal of Climate},
volume = {32},
number = {10},
pages = {2917-2949},
year = {2019},
doi = {10.1175/JCLI-D-18-0796.1},
URL = {https://doi.org/10.1175/JCLI-D-18-0796.1},
eprint = {https://doi.org/10.1175/JCLI-D-18-0796.1}
}"
"@

$bibJAS0233 = @"
 "@article{doi:10.1175/JAS-D-13-0233.1,
author = {Park, Sungsu},
title = {A Unified Convection Scheme (UNICON). Part I: Formulation},
journal = {Journal of the Atmospheric Sciences},
volume = {71},
number = {11},
pages = {3902-3930},
year = {2014},
doi = {10.1175/JAS-D-13-0233.1},
URL = {https://doi.org/10.1175/JAS-D-13-0233.1},
eprint = {https://doi.org/10.1175/JAS-D-13-0233.1}
}"
"@

$bibJAS0234 = @"
 "@article{doi:10.1175/JAS-D-13-0234.1,
author = {Park, Sungsu},
title = {A Unified Convection Scheme (UNICON). Part II: Simulation},
journal = {Journal of the Atmospheric Sciences},
volume = {71},
number = {11},
pages = {3931-3973},
year = {2014},
doi = {10.1175/JAS-D-13-0234.1},
URL = {https://doi.org/10.1175/JAS-D-13-0234.1},
eprint = {https://doi.org/10.1175/JAS-D-13-0234.1}
}"
"@

$bibMS000877 = @"
 "@article{doi:10.1002/2016MS000877,
author = {Park, Sungsu and Baek, Eun-Hyuk and Kim, Baek-Min and Kim, Seong-Joong},
title = {Impact of detrained cumulus on climate simulated by the Community Atmosphere Model Version 5 with a unified convection scheme},
journal = {Journal of Advances in Modeling Earth Systems},
volume = {9},
number = {2},
pages = {1399-1411},
keywords = {detrained cumulus, unified convection, low-level cloud, parameterization},
doi = {10.1002/2016MS000877},
url = {https://agupubs.onlinelibrary.wiley.com/doi/abs/10.1002/2016MS000877},
eprint = {https://agupubs.onlinelibrary.wiley.com/doi/pdf/10.1002/2016MS000877},
year = {2017}
}"
"@

$bibJCLI = $bibJCLI.TrimStart()
$bibJAS0233 = $bibJAS0233.TrimStart()
$bibJAS0234 = $bibJAS0234.TrimStart()
$bibMS000877 = $bibMS000877.TrimStart()

# ---------------------------------------------------------------------
# Row 3 : SAM0UNICON_2019 / JCLI-D-18-0796.1
# ---------------------------------------------------------------------
$ws.Range("A3").Value = "SAM0UNICON_2019"

$ws.Range("B3").Value = "https://doi.org/10.1175/JCLI-D-18-0796.1"
$ws.Hyperlinks.Add($ws.Range("B3"), "https://doi.org/10.1175/JCLI-D-18-0796.1")

$ws.Range("C3").Value = $bibJCLI

$ws.Range("D3").Borders.LineStyle = -4142
$ws.Range("D3").HorizontalAlignment = 1
$ws.Range("D3").Value = "https://journals.ametsoc.org/doi/full/10.1175/JCLI-D-18-0796.1"
$ws.Hyperlinks.Add($ws.Range("D3"), "https://journals.ametsoc.org/doi/full/10.1175/JCLI-D-18-0796.1")

$ws.Range("E3").Value = "Paper with full description of CMIP6 simulation results"

$ws.Rows.Item(3).RowHeight = 165.75

# ---------------------------------------------------------------------
# Row 4 : UNICON_2014a / JAS-D-13-0233.1
# ---------------------------------------------------------------------
$ws.Range("A4").Value = "UNICON_2014a"

$ws.Range("B4").Value = "https://doi.org/10.1175/JAS-D-13-0233.1"
$ws.Hyperlinks.Add($ws.Range("B4"), "https://doi.org/10.1175/JAS-D-13-0233.1")

$ws.Range("C4").Value = $bibJAS0233

$ws.Range("D4").Borders.LineStyle = -4142
$ws.Range("D4").HorizontalAlignment = 1
$ws.Range("D4").Value = "https://journals.ametsoc.org/doi/full/10.1175/JAS-D-13-0233.1"
$ws.Hyperlinks.Add($ws.Range("D4"), "https://journals.ametsoc.org/doi/full/10.1175/JAS-D-13-0233.1")

$ws.Range("E4").Value = "Paper about UNICON"

$ws.Rows.Item(4).RowHeight = 153

# ---------------------------------------------------------------------
# Row 5 : UNICON_2014b / JAS-D-13-0234.1
# ---------------------------------------------------------------------
$ws.Range("A5").Value = "UNICON_2014b"

$ws.Range("B5").Value = "https://doi.org/10.1175/JAS-D-13-0234.1"
$ws.Hyperlinks.Add($ws.Range("B5"), "https://doi.org/10.1175/JAS-D-13-0234.1")

$ws.Range("C5").Value = $bibJAS0234

$ws.Range("D5").Borders.LineStyle = -4142
$ws.Range("D5").HorizontalAlignment = 1
$ws.Range("D5").Value = "https://journals.ametsoc.org/doi/full/10.1175/JAS-D-13-0234.1"
$ws.Hyperlinks.Add($ws.Range("D5"), "https://journals.ametsoc.org/doi/full/10.1175/JAS-D-13-0234.1")

$ws.Rows.Item(5).RowHeight = 153

# ---------------------------------------------------------------------
# Row 6 : UNICON_2017 / 2016MS000877
# ---------------------------------------------------------------------
$ws.Range("A6").Value = "UNICON_2017"

$ws.Range("B6").Borders.LineStyle = -4142
$ws.Range("B6").HorizontalAlignment = 1
$ws.Range("B6").VerticalAlignment = -4108
$ws.Range("B6").Value = "https://doi.org/10.1002/2016MS000877"
$ws.Hyperlinks.Add($ws.Range("B6"), "https://doi.org/10.1002/2016MS000877")

$ws.Range("C6").Value = $bibMS000877

$ws.Range("D6").Borders.LineStyle = -4142
$ws.Range("D6").HorizontalAlignment = 1
$ws.Range("D6").Value = "https://agupubs.onlinelibrary.wiley.com/doi/full/10.1002/2016MS000877"
$ws.Hyperlinks.Add($ws.Range("D6"), "https://agupubs.onlinelibrary.wiley.com/doi/full/10.1002/2016MS000877")

$ws.Rows.Item(6).RowHeight = 178.5

# ---------------------------------------------------------------------
# Row 7 : carry B6's (non-bordered / vertically centered hyperlink) look
# down into B7, without any content or hyperlink of its own.
# ---------------------------------------------------------------------
$ws.Range("B6").Copy()
$ws.Range("B7").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# Make "Citations" the active sheet/tab, as in the edited workbook.
# ---------------------------------------------------------------------
$ws.Activate()
$ws.Range("A3").Select()

Write-Host "Citations sheet updated"
